$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "51.156.63"
$ws.Range("E2").Value = "  -15.83%  "

$ws.Range("D3").Value = "2.272.64"
$ws.Range("E3").Value = "  -21.81%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.27%  "

Set-TextValue $ws.Range("D5") "425.46"
$ws.Range("E5").Value = "  -19.15%  "

Set-TextValue $ws.Range("D6") "118.41"
$ws.Range("E6").Value = "  -17.94%  "

Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.13%  "

Set-TextValue $ws.Range("D8") "0.455"
$ws.Range("E8").Value = "  -16.83%  "

$ws.Range("D9").Value = "2.286.84"
$ws.Range("E9").Value = "  -21.52%  "

Set-TextValue $ws.Range("D10") "5.04"
$ws.Range("E10").Value = "  -17.93%  "

Set-TextValue $ws.Range("D11") "0.0871"
$ws.Range("E11").Value = "  -18.76%  "

Set-TextValue $ws.Range("D12") "0.298"
$ws.Range("E12").Value = "  -16.75%  "

$ws.Range("E13").Value = "  -7.08%  "

$ws.Range("D14").Value = "2.654.84"
$ws.Range("E14").Value = "  -22.18%  "

$ws.Range("D15").Value = "51.085.80"
$ws.Range("E15").Value = "  -15.97%  "

Set-TextValue $ws.Range("D16") "18.42"
$ws.Range("E16").Value = "  -18.27%  "

Set-TextValue $ws.Range("D17") "0.0000116"
$ws.Range("E17").Value = "  -18.05%  "

$ws.Range("D18").Value = "2.275.47"
$ws.Range("E18").Value = "  -21.72%  "

Set-TextValue $ws.Range("D19") "3.98"
$ws.Range("E19").Value = "  -18.67%  "

Set-TextValue $ws.Range("D20") "292.08"
$ws.Range("E20").Value = "  -17.31%  "

Set-TextValue $ws.Range("D21") "1.00"
$ws.Range("E21").Value = "  +0.17%  "

Set-TextValue $ws.Range("D22") "5.67"
$ws.Range("E22").Value = "  -0.43%  "

Set-TextValue $ws.Range("D23") "8.62"
$ws.Range("E23").Value = "  -25.38%  "

Set-TextValue $ws.Range("D24") "5.08"
$ws.Range("E24").Value = "  -21.93%  "

Set-TextValue $ws.Range("D25") "0.988"
$ws.Range("E25").Value = "  -1.40%  "

Set-TextValue $ws.Range("D26") "52.43"
$ws.Range("E26").Value = "  -19.08%  "

Set-TextValue $ws.Range("D27") "0.365"
$ws.Range("E27").Value = "  -18.92%  "

$ws.Range("D28").Value = "2.309.50"
$ws.Range("E28").Value = "  -23.86%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D29") "0.134"
$ws.Range("E29").Value = "  -25.08%  "

$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D30") "0.997"
$ws.Range("E30").Value = "  -0.26%  "

Set-TextValue $ws.Range("D31") "6.64"
$ws.Range("E31").Value = "  -14.96%  "

$ws.Range("D32").Value = "0.0₃0644"
$ws.Range("E32").Value = "  -25.92%  "

Set-TextValue $ws.Range("D33") "142.82"
$ws.Range("E33").Value = "  -6.82%  "

Set-TextValue $ws.Range("D34") "16.58"
$ws.Range("E34").Value = "  -15.54%  "

Set-TextValue $ws.Range("D35") "1.29"
$ws.Range("E35").Value = "  -23.20%  "

Set-TextValue $ws.Range("D36") "4.70"
$ws.Range("E36").Value = "  -15.91%  "

Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.30%  "

Set-TextValue $ws.Range("D38") "3.24"
$ws.Range("E38").Value = "  -26.34%  "

Set-TextValue $ws.Range("D39") "0.970"
$ws.Range("E39").Value = "  -19.10%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D40") "31.94"
$ws.Range("E40").Value = "  -14.98%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D41") "0.737"
$ws.Range("E41").Value = "  -25.96%  "

$ws.Range("E42").Value = "  -2.20%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.551"
$ws.Range("E43").Value = "  -15.57%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D44") "3.11"
$ws.Range("E44").Value = "  -15.82%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D45") "0.0491"
$ws.Range("E45").Value = "  -15.75%  "

$ws.Range("D46").Value = "1.867.61"
$ws.Range("E46").Value = "  -18.51%  "

Set-TextValue $ws.Range("D47") "1.12"
$ws.Range("E47").Value = "  -23.70%  "

Set-TextValue $ws.Range("D48") "0.0199"
$ws.Range("E48").Value = "  -16.04%  "

Set-TextValue $ws.Range("D49") "0.0795"
$ws.Range("E49").Value = "  -13.05%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "15.51"
$ws.Range("E50").Value = "  -23.88%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D51") "3.86"
$ws.Range("E51").Value = "  -21.79%  "
